# BIS-1002: removed "Internal Assignment" column from export.
# The "Internal Assignment" header lived in O4 with the per-row values in
# O5:O7. Clearing those cells' contents removes the column's data while
# leaving the (now-empty) styled cells in place, and drops the
# now-unreferenced "Internal Assignment" shared string on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4:O7").ClearContents()
